$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overall": A2 (477) was a numeric cell -> convert to text "477"
# ---------------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
$wsOverall.Range("A2").Value = "'477"

# ---------------------------------------------------------------------------
# Sheet "County": column B (filer counts) numeric -> text, for every county
# row (2-41). Leading apostrophe forces a literal text value instead of
# Excel's automatic numeric re-parsing.
# ---------------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")
$countyCounts = @{
    2  = "161"
    3  = "1"
    4  = "21"
    5  = "2"
    6  = "12"
    7  = "25"
    8  = "1"
    9  = "17"
    10 = "22"
    11 = "4"
    12 = "1"
    13 = "1"
    14 = "29"
    15 = "3"
    16 = "1"
    17 = "3"
    18 = "2"
    19 = "4"
    20 = "2"
    21 = "5"
    22 = "3"
    23 = "2"
    24 = "4"
    25 = "1"
    26 = "3"
    27 = "41"
    28 = "27"
    29 = "10"
    30 = "1"
    31 = "2"
    32 = "2"
    33 = "12"
    34 = "1"
    35 = "2"
    36 = "4"
    37 = "3"
    38 = "12"
    39 = "22"
    40 = "5"
    41 = "3"
}
foreach ($r in $countyCounts.Keys) {
    $wsCounty.Range("B$r").Value = "'" + $countyCounts[$r]
}

# Rows 42 (Power County) and 43 (Bear Lake County) previously held all-zero
# placeholder data ("0" in every column). They now carry properly formatted
# placeholder text matching the rest of the sheet's percentage / currency
# formatting.
$wsCounty.Range("B42").Value = "'0.00%"
$wsCounty.Range("C42").Value = "'`$0"
$wsCounty.Range("D42").Value = "'0.00%"
$wsCounty.Range("E42").Value = "'0.00%"
$wsCounty.Range("F42").Value = "'0.00%"

$wsCounty.Range("B43").Value = "'0.00%"
$wsCounty.Range("C43").Value = "'`$0"
$wsCounty.Range("D43").Value = "'0.00%"
$wsCounty.Range("E43").Value = "'0.00%"
$wsCounty.Range("F43").Value = "'0.00%"

# New row 44: a statewide "Total" row appended below the county list
# (extends the sheet dimension from A1:F43 to A1:F44).
$wsCounty.Range("A44").Value = "Total"
$wsCounty.Range("B44").Value = "'477"
$wsCounty.Range("C44").Value = "'`$643,142,727"
$wsCounty.Range("D44").Value = "'12.56%"
$wsCounty.Range("E44").Value = "'-11.31%"
$wsCounty.Range("F44").Value = "'62.05%"

# ---------------------------------------------------------------------------
# Sheet "Congressional District": column B numeric -> text (rows 2-4)
# ---------------------------------------------------------------------------
$wsCd = $wb.Worksheets.Item("Congressional District")
$cdCounts = @{
    2 = "194"
    3 = "283"
    4 = "477"
}
foreach ($r in $cdCounts.Keys) {
    $wsCd.Range("B$r").Value = "'" + $cdCounts[$r]
}

# ---------------------------------------------------------------------------
# Sheet "Size": column B numeric -> text (rows 2-8)
# ---------------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
$sizeCounts = @{
    2 = "183"
    3 = "138"
    4 = "79"
    5 = "27"
    6 = "37"
    7 = "13"
    8 = "477"
}
foreach ($r in $sizeCounts.Keys) {
    $wsSize.Range("B$r").Value = "'" + $sizeCounts[$r]
}

# ---------------------------------------------------------------------------
# Sheet "Subsector": column B numeric -> text (rows 2-13)
# ---------------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subsector")
$subCounts = @{
    2  = "36"
    3  = "93"
    4  = "41"
    5  = "36"
    6  = "8"
    7  = "124"
    8  = "6"
    9  = "31"
    10 = "11"
    11 = "88"
    12 = "3"
    13 = "477"
}
foreach ($r in $subCounts.Keys) {
    $wsSub.Range("B$r").Value = "'" + $subCounts[$r]
}
